$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the 5 runs of the "The call went through..." paragraph into one.
#    A Find/Replace (even a no-op one) that touches the first run causes the
#    engine to coalesce adjacent same-formatted runs in the paragraph.
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("Gary agreed", $true, $false, $false, $false, $false, $true, 1, $false, "Gary agreed", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Fix the "resently" -> "presently" typo, then split that run into three
#    runs: "...assembler would ", "presently", " accept. ...".
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("resently", $true, $false, $false, $false, $false, $true, 1, $false, "presently", 2) | Out-Null

$r2b = $d.Content
$r2b.Find.Execute("presently", $true) | Out-Null
$pStart = $r2b.Start
$pEnd = $r2b.End

$d.Bookmarks.Add("_TmpSplitA", $d.Range($pStart, $pStart))
$d.Bookmarks.Add("_TmpSplitB", $d.Range($pEnd, $pEnd))
$d.Bookmarks("_TmpSplitA").Delete()
$d.Bookmarks("_TmpSplitB").Delete()

# ---------------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from the end of the document to the middle
#    of the word "computer" ("low-level compute" | "r, and the I/O...").
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("low-level computer", $true) | Out-Null
$splitPoint = $r3.Start + 17
$d.Bookmarks.Add("_GoBack", $d.Range($splitPoint, $splitPoint))

# ---------------------------------------------------------------------------
# 4) Fix the "Intellec8 mode 80" -> "Intellec8 mod 80" typo.
# ---------------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("Intellec8 mode 80", $true, $false, $false, $false, $false, $true, 1, $false, "Intellec8 mod 80", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Merge the 4 runs of the "By comparison, ..." paragraph into one.
# ---------------------------------------------------------------------------
$r5 = $d.Content
$r5.Find.Execute("By comparison, t", $true, $false, $false, $false, $false, $true, 1, $false, "By comparison, t", 2) | Out-Null
